$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 974, shifting existing data
# (rows 974:1045) down to (976:1047).
$ws.Rows("974:975").Insert()

# Populate the two newly inserted rows with the new weekly price record
# (same shape/columns as the surrounding rows, date 45021).

# Row 974 - "Primera" quality
$ws.Cells.Item(974, 1).Value = 3
$ws.Cells.Item(974, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(974, 3).Value = "Coquimbo"
$ws.Cells.Item(974, 4).Value = 45021
$ws.Cells.Item(974, 5).Value = 5
$ws.Cells.Item(974, 6).Value = 100114014
$ws.Cells.Item(974, 7).Value = "Betarraga"
$ws.Cells.Item(974, 8).Value = "Sin especificar"
$ws.Cells.Item(974, 9).Value = "Primera"
$ws.Cells.Item(974, 10).Value = 3100
$ws.Cells.Item(974, 11).Value = 750
$ws.Cells.Item(974, 12).Value = 800
$ws.Cells.Item(974, 13).Value = 776
$ws.Cells.Item(974, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(974, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(974, 16).Value = 194
$ws.Cells.Item(974, 17).Value = 4
$ws.Cells.Item(974, 18).Value = "Hortaliza"

# Row 975 - "Segunda" quality
$ws.Cells.Item(975, 1).Value = 3
$ws.Cells.Item(975, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(975, 3).Value = "Coquimbo"
$ws.Cells.Item(975, 4).Value = 45021
$ws.Cells.Item(975, 5).Value = 5
$ws.Cells.Item(975, 6).Value = 100114014
$ws.Cells.Item(975, 7).Value = "Betarraga"
$ws.Cells.Item(975, 8).Value = "Sin especificar"
$ws.Cells.Item(975, 9).Value = "Segunda"
$ws.Cells.Item(975, 10).Value = 1600
$ws.Cells.Item(975, 11).Value = 600
$ws.Cells.Item(975, 12).Value = 600
$ws.Cells.Item(975, 13).Value = 600
$ws.Cells.Item(975, 14).Value = "$/paquete 4 unidades"
$ws.Cells.Item(975, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(975, 16).Value = 150
$ws.Cells.Item(975, 17).Value = 4
$ws.Cells.Item(975, 18).Value = "Hortaliza"
